$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.275.85"
$ws.Range("E2").Value = "  +2.38%  "
$ws.Range("D3").Value = "3.386.59"
$ws.Range("E3").Value = "  +1.72%  "
$ws.Range("D4").Value = "'" + "1.00"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'" + "586.29"
$ws.Range("E5").Value = "  +1.07%  "
$ws.Range("D6").Value = "'" + "180.20"
$ws.Range("E6").Value = "  +2.57%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("E8").Value = "  +1.49%  "
$ws.Range("E9").Value = "  +8.52%  "
$ws.Range("E10").Value = "  +2.19%  "
$ws.Range("D11").Value = "'" + "48.51"
$ws.Range("E11").Value = "  +3.48%  "
$ws.Range("D12").Value = "'" + "0.0000283"
$ws.Range("E12").Value = "  +4.21%  "
$ws.Range("D13").Value = "'" + "678.71"
$ws.Range("E13").Value = "  -1.82%  "
$ws.Range("D14").Value = "'" + "8.65"
$ws.Range("E14").Value = "  +3.05%  "
$ws.Range("D15").Value = "3.927.25"
$ws.Range("E15").Value = "  +1.43%  "
$ws.Range("D16").Value = "69.360.22"
$ws.Range("E16").Value = "  +2.51%  "
$ws.Range("B17").Value = "TRON"
$ws.Range("C17").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D17").Value = "'" + "0.120"
$ws.Range("E17").Value = "  +1.44%  "
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "3.381.11"
$ws.Range("E18").Value = "  +1.55%  "
$ws.Range("D19").Value = "'" + "17.71"
$ws.Range("E19").Value = "  +0.67%  "
$ws.Range("E20").Value = "  +2.32%  "
$ws.Range("D21").Value = "'" + "0.906"
$ws.Range("E21").Value = "  +1.56%  "
$ws.Range("D22").Value = "'" + "5.40"
$ws.Range("E22").Value = "  -2.30%  "
$ws.Range("D23").Value = "'" + "17.18"
$ws.Range("E23").Value = "  +1.87%  "
$ws.Range("D24").Value = "'" + "103.35"
$ws.Range("E24").Value = "  +2.36%  "
$ws.Range("E25").Value = "  +0.53%  "
$ws.Range("D26").Value = "'" + "2.72"
$ws.Range("E26").Value = "  +1.71%  "
$ws.Range("E27").Value = "  +2.96%  "
$ws.Range("D28").Value = "'" + "33.85"
$ws.Range("E28").Value = "  +2.21%  "
$ws.Range("D29").Value = "'" + "8.80"
$ws.Range("E29").Value = "  +3.05%  "
$ws.Range("D30").Value = "'" + "6.95"
$ws.Range("E30").Value = "  -1.60%  "
$ws.Range("D31").Value = "'" + "11.15"
$ws.Range("E31").Value = "  +1.41%  "
$ws.Range("B32").Value = "Bittensor"
$ws.Range("C32").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D32").Value = "'" + "557.57"
$ws.Range("E32").Value = "  -2.26%  "
$ws.Range("B33").Value = "dogwifhat"
$ws.Range("C33").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D33").Value = "'" + "3.61"
$ws.Range("E33").Value = "  +10.46%  "
$ws.Range("D34").Value = "'" + "0.106"
$ws.Range("E34").Value = "  +1.21%  "
$ws.Range("D35").Value = "'" + "58.53"
$ws.Range("E35").Value = "  +2.20%  "
$ws.Range("E36").Value = "  +0.13%  "
$ws.Range("D37").Value = "3.667.82"
$ws.Range("E37").Value = "  -1.17%  "
$ws.Range("E38").Value = "  +5.47%  "
$ws.Range("D39").Value = "'" + "35.55"
$ws.Range("E39").Value = "  +0.66%  "
$ws.Range("D40").Value = "0.0₃0719"
$ws.Range("E40").Value = "  +7.63%  "
$ws.Range("D41").Value = "'" + "3.27"
$ws.Range("E41").Value = "  +3.78%  "
$ws.Range("D42").Value = "'" + "2.68"
$ws.Range("E42").Value = "  +2.89%  "
$ws.Range("D43").Value = "'" + "0.339"
$ws.Range("E43").Value = "  +1.43%  "
$ws.Range("D44").Value = "'" + "0.0424"
$ws.Range("E44").Value = "  +4.17%  "
$ws.Range("B45").Value = "ThetaToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D45").Value = "'" + "2.69"
$ws.Range("E45").Value = "  +1.66%  "
$ws.Range("B46").Value = "Stellar"
$ws.Range("C46").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D46").Value = "'" + "0.130"
$ws.Range("E46").Value = "  +1.25%  "
$ws.Range("B47").Value = "Mantle"
$ws.Range("C47").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D47").Value = "'" + "1.40"
$ws.Range("E47").Value = "  +5.03%  "
$ws.Range("B48").Value = "FirstDigitalUSD"
$ws.Range("C48").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D48").Value = "'" + "1.00"
$ws.Range("E48").Value = "  -0.08%  "
$ws.Range("B49").Value = "Monero"
$ws.Range("C49").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D49").Value = "'" + "134.13"
$ws.Range("E49").Value = "  +1.60%  "
$ws.Range("B50").Value = "CoreDAO"
$ws.Range("C50").Value = "https://coinranking.com/coin/HFvoXUQh4+coredao-core"
$ws.Range("D50").Value = "'" + "2.65"
$ws.Range("E50").Value = "  +3.31%  "
$ws.Range("B51").Value = "THORChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D51").Value = "'" + "7.54"
$ws.Range("E51").Value = "  +3.07%  "
